# Apply content rewrite for the two paragraphs (merging/splitting runs,
# dropping explicit color/underline overrides, adding hanging-indent
# paragraph formatting, and updating the text per the commit diff).
$d = $word.ActiveDocument

$xml = @'
<?xml version="1.0" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:wpc="http://schemas.microsoft.com/office/word/2010/wordprocessingCanvas" xmlns:mc="http://schemas.openxmlformats.org/markup-compatibility/2006" xmlns:o="urn:schemas-microsoft-com:office:office" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:m="http://schemas.openxmlformats.org/officeDocument/2006/math" xmlns:v="urn:schemas-microsoft-com:vml" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" xmlns:w10="urn:schemas-microsoft-com:office:word" xmlns:w15="http://schemas.microsoft.com/office/word/2012/wordml" xmlns:wpg="http://schemas.microsoft.com/office/word/2010/wordprocessingGroup" xmlns:wpi="http://schemas.microsoft.com/office/word/2010/wordprocessingInk" xmlns:wne="http://schemas.microsoft.com/office/word/2006/wordml" xmlns:wps="http://schemas.microsoft.com/office/word/2010/wordprocessingShape" mc:Ignorable="w14 w15 wp14"><w:body><w:p>
      <w:pPr>
        <w:ind w:left="100" w:hanging="100" w:hangingChars="50"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t xml:space="preserve">  +Cách chuyển từ Bit sang Byte ta lấy Bit chia cho 8 </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:br w:type="textWrapping"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t xml:space="preserve">+  1 byte = 8bit </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:br w:type="textWrapping"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t>=&gt; 1024 bit = 1024/8 byte</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="100" w:hanging="100" w:hangingChars="50"/>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t>=&gt; 1024bit = 128 byte</w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:br w:type="textWrapping"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t xml:space="preserve">lý giải 1 byte bằng 8 bit : + trong nhu cầu mã hóa kí tự 7 bit đủ để biểu diễn tất cả các ký tư tiếng ah nhưng với 8bit sẽ linh hoạt hơn và là một con số tốt hơn để xử lý dữ diệu </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:br w:type="textWrapping"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t xml:space="preserve">8bit sẽ tạo ra 2^8 kí tự đủ để mã hóa một ký tự trong bằng mã ASCII  bao gồm chữ cái , số và dữ liệu </w:t>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:br w:type="textWrapping"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:br w:type="textWrapping"/>
      </w:r>
      <w:r>
        <w:rPr>
          <w:rFonts w:hint="default"/>
          <w:lang w:val="vi-VN"/>
        </w:rPr>
        <w:t xml:space="preserve">                                            </w:t>
      </w:r>
      <w:bookmarkStart w:id="0" w:name="_GoBack"/>
      <w:bookmarkEnd w:id="0"/>
    </w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$d.Content.InsertXML($xml)

# Remove the now-unused built-in "Hyperlink" character style.
$d.Styles("Hyperlink").Delete()
